$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(5, 5).Value = 0.0000006  # E5
$ws.Cells.Item(10, 2).Value = 0.02109048192771084  # B10
$ws.Cells.Item(11, 2).Value = 0.00534433734939759  # B11
$ws.Cells.Item(12, 2).Value = 0.03012048192771084  # B12
$ws.Cells.Item(13, 2).Value = 0.03012048192771084  # B13
$ws.Cells.Item(17, 2).Value = 0.02228506024096386  # B17
$ws.Cells.Item(17, 5).Value = 0.0030157  # E17
$ws.Cells.Item(18, 2).Value = 0.05067674698795181  # B18
$ws.Cells.Item(18, 5).Value = 0.0164768  # E18
$ws.Cells.Item(19, 2).Value = 0.06156084337349398  # B19
$ws.Cells.Item(19, 5).Value = 0.0384108  # E19
$ws.Cells.Item(20, 2).Value = 0.1205855421686747  # B20
$ws.Cells.Item(20, 5).Value = 0.0222614  # E20
$ws.Cells.Item(21, 2).Value = 0.03525915662650603  # B21
$ws.Cells.Item(21, 5).Value = 0.0000001  # E21
$ws.Cells.Item(21, 7).Value = 0.17999738  # G21
$ws.Cells.Item(22, 5).Value = 0.0000001  # E22
$ws.Cells.Item(24, 2).Value = 0  # B24
$ws.Cells.Item(24, 4).Value = 0.25  # D24
$ws.Cells.Item(24, 6).Value = 0.2123866  # F24
$ws.Cells.Item(25, 2).Value = 0  # B25
$ws.Cells.Item(25, 6).Value = 0.07339333333333334  # F25
$ws.Cells.Item(26, 2).Value = 0.02181060240963855  # B26
$ws.Cells.Item(26, 6).Value = 0.05053333333333333  # F26
$ws.Cells.Item(27, 2).Value = 0.01433385542168675  # B27
$ws.Cells.Item(27, 5).Value = 0  # E27
$ws.Cells.Item(27, 6).Value = 0.04666666666666667  # F27
$ws.Cells.Item(28, 2).Value = 0  # B28
$ws.Cells.Item(28, 5).Value = 0  # E28
$ws.Cells.Item(28, 7).Value = 0.27883965  # G28
$ws.Cells.Item(29, 2).Value = 0.1111146987951807  # B29
$ws.Cells.Item(29, 5).Value = 0  # E29
$ws.Cells.Item(30, 2).Value = 0.1175892771084337  # B30
$ws.Cells.Item(30, 5).Value = 0  # E30
$ws.Cells.Item(31, 2).Value = 0.01604361445783132  # B31
$ws.Cells.Item(31, 5).Value = 0  # E31
$ws.Cells.Item(37, 2).Value = 0.005695542168674699  # B37
$ws.Cells.Item(38, 2).Value = 0.03520373493975903  # B38
$ws.Cells.Item(39, 2).Value = 0.0193389156626506  # B39
$ws.Cells.Item(43, 2).Value = 0.02456301204819277  # B43
$ws.Cells.Item(45, 2).Value = 0.05421686746987952  # B45
$ws.Cells.Item(45, 5).Value = 0.007922200000000001  # E45
$ws.Cells.Item(46, 2).Value = 0.04735578313253012  # B46
$ws.Cells.Item(46, 5).Value = 0.0488011  # E46
$ws.Cells.Item(47, 2).Value = 0.01251759036144578  # B47
$ws.Cells.Item(47, 5).Value = 0.0533715  # E47
$ws.Cells.Item(48, 2).Value = 0.01685156626506024  # B48
$ws.Cells.Item(48, 5).Value = 0.0401201  # E48
$ws.Cells.Item(49, 2).Value = 0.0132689156626506  # B49
$ws.Cells.Item(49, 4).Value = 0.375  # D49
$ws.Cells.Item(49, 5).Value = 0.0000002  # E49
$ws.Cells.Item(50, 2).Value = 0.006185542168674699  # B50
$ws.Cells.Item(51, 2).Value = 0.02420927710843373  # B51
$ws.Cells.Item(51, 5).Value = 0  # E51
$ws.Cells.Item(52, 2).Value = 0.02382204819277109  # B52
$ws.Cells.Item(52, 5).Value = 0  # E52
$ws.Cells.Item(53, 2).Value = 0.01266614457831325  # B53
$ws.Cells.Item(53, 5).Value = 0  # E53
$ws.Cells.Item(54, 2).Value = 0.0125110843373494  # B54
$ws.Cells.Item(54, 5).Value = 0  # E54
$ws.Cells.Item(55, 2).Value = 0.0228678313253012  # B55
$ws.Cells.Item(55, 5).Value = 0  # E55
$ws.Cells.Item(56, 2).Value = 0.03629228915662651  # B56
$ws.Cells.Item(56, 5).Value = 0  # E56
$ws.Cells.Item(57, 2).Value = 0.01133156626506024  # B57
$ws.Cells.Item(58, 2).Value = 0.02426650602409638  # B58
$ws.Cells.Item(59, 2).Value = 0.01683686748192771  # B59
$ws.Cells.Item(60, 2).Value = 0.007806024096385542  # B60
$ws.Cells.Item(68, 2).Value = 0.02049759036144578  # B68
$ws.Cells.Item(69, 2).Value = 0.003598795180722892  # B69
$ws.Cells.Item(70, 7).Value = 0.195854  # G70
$ws.Cells.Item(71, 2).Value = 0.01885734939759036  # B71
$ws.Cells.Item(71, 7).Value = 0.07990800000000001  # G71
$ws.Cells.Item(72, 2).Value = 0.1055784337349398  # B72
$ws.Cells.Item(73, 2).Value = 0.08380771083132529  # B73
$ws.Cells.Item(73, 5).Value = 0.048495  # E73
$ws.Cells.Item(74, 2).Value = 0.1156054216987952  # B74
$ws.Cells.Item(74, 5).Value = 0.067893  # E74
$ws.Cells.Item(75, 2).Value = 0.03385373493975904  # B75
$ws.Cells.Item(75, 4).Value = 0.006410256410000001  # D75
$ws.Cells.Item(75, 5).Value = 0.067893  # E75
$ws.Cells.Item(76, 4).Value = 0.044871794875  # D76
$ws.Cells.Item(76, 5).Value = 0.067893  # E76
$ws.Cells.Item(77, 2).Value = 0  # B77
$ws.Cells.Item(77, 3).Value = 0  # C77
$ws.Cells.Item(77, 4).Value = 0.044871794875  # D77
$ws.Cells.Item(77, 5).Value = 0.067893  # E77
$ws.Cells.Item(78, 2).Value = 0  # B78
$ws.Cells.Item(78, 3).Value = 0  # C78
$ws.Cells.Item(78, 4).Value = 0.044871794875  # D78
$ws.Cells.Item(78, 5).Value = 0.067893  # E78
$ws.Cells.Item(78, 7).Value = 0.383153  # G78
$ws.Cells.Item(79, 2).Value = 0.03280939759036144  # B79
$ws.Cells.Item(79, 3).Value = 0  # C79
$ws.Cells.Item(79, 4).Value = 0.044871794875  # D79
$ws.Cells.Item(79, 5).Value = 0.067893  # E79
$ws.Cells.Item(80, 2).Value = 0.04433048192771084  # B80
$ws.Cells.Item(80, 3).Value = 0  # C80
$ws.Cells.Item(80, 4).Value = 0.044871794875  # D80
$ws.Cells.Item(80, 5).Value = 0.067893  # E80
$ws.Cells.Item(80, 8).Value = 0  # H80
$ws.Cells.Item(81, 2).Value = 0.04220819277108434  # B81
$ws.Cells.Item(81, 3).Value = 0  # C81
$ws.Cells.Item(81, 4).Value = 0.044871794875  # D81
$ws.Cells.Item(81, 5).Value = 0.067893  # E81
$ws.Cells.Item(82, 2).Value = 0.02282060240963855  # B82
$ws.Cells.Item(82, 3).Value = 0  # C82
$ws.Cells.Item(82, 4).Value = 0.044871794875  # D82
$ws.Cells.Item(82, 5).Value = 0.067893  # E82
$ws.Cells.Item(83, 2).Value = 0  # B83
$ws.Cells.Item(83, 3).Value = 0  # C83
$ws.Cells.Item(83, 4).Value = 0.01282051281666667  # D83
$ws.Cells.Item(83, 5).Value = 0.019398  # E83
$ws.Cells.Item(84, 7).Value = 0.20124  # G84
$ws.Cells.Item(86, 2).Value = 0.01324614457831325  # B86
$ws.Cells.Item(87, 2).Value = 0.01676481927710843  # B87
$ws.Cells.Item(88, 2).Value = 0.01280807228915663  # B88
$ws.Cells.Item(89, 2).Value = 0.01260265060240964  # B89
$ws.Cells.Item(90, 2).Value = 0.01435566265060241  # B90
$ws.Cells.Item(91, 2).Value = 0.0109455421686747  # B91
$ws.Cells.Item(92, 2).Value = 0.009638554216867471  # B92
$ws.Cells.Item(93, 2).Value = 0.009638554216867471  # B93
$ws.Cells.Item(101, 6).Value = 0  # F101
$ws.Cells.Item(124, 7).Value = 0.075  # G124
$ws.Cells.Item(125, 2).Value = 0.08781590360240964  # B125
$ws.Cells.Item(125, 4).Value = 0.000000040516486725  # D125
$ws.Cells.Item(125, 5).Value = 0.100597266061424  # E125
$ws.Cells.Item(125, 7).Value = 0.141738  # G125
$ws.Cells.Item(126, 2).Value = 0.1114400000481928  # B126
$ws.Cells.Item(126, 5).Value = 0.16220282044  # E126
$ws.Cells.Item(127, 2).Value = 0.1191290361445783  # B127
$ws.Cells.Item(127, 5).Value = 0.12399304737  # E127
$ws.Cells.Item(128, 2).Value = 0.1227467469879518  # B128
$ws.Cells.Item(128, 5).Value = 0.05983448736  # E128
$ws.Cells.Item(129, 2).Value = 0.3247060240963855  # B129
$ws.Cells.Item(129, 5).Value = 0.0432356567  # E129
$ws.Cells.Item(129, 7).Value = 0.08  # G129
$ws.Cells.Item(130, 2).Value = 0.2872278313253012  # B130
$ws.Cells.Item(130, 5).Value = 0.08690113362  # E130
$ws.Cells.Item(130, 7).Value = 0.348549  # G130
$ws.Cells.Item(131, 2).Value = 0.06737120481927711  # B131
$ws.Cells.Item(131, 5).Value = 0.0455531954  # E131
$ws.Cells.Item(132, 2).Value = 0  # B132
$ws.Cells.Item(134, 2).Value = 0.198054578313253  # B134
$ws.Cells.Item(135, 2).Value = 0.1305381927710843  # B135
$ws.Cells.Item(137, 7).Value = 0.22815  # G137
$ws.Cells.Item(149, 2).Value = 0  # B149
$ws.Cells.Item(150, 2).Value = 0.1029344578313253  # B150
$ws.Cells.Item(151, 2).Value = 0.1701113253012048  # B151
$ws.Cells.Item(152, 2).Value = 0.1910298795180723  # B152
$ws.Cells.Item(153, 2).Value = 0.1264832530120482  # B153
$ws.Cells.Item(153, 5).Value = 0  # E153
$ws.Cells.Item(154, 2).Value = 0.04967542168674699  # B154
$ws.Cells.Item(154, 5).Value = 0  # E154
$ws.Cells.Item(155, 2).Value = 0.02241481927710844  # B155
$ws.Cells.Item(155, 5).Value = 0  # E155
$ws.Cells.Item(155, 8).Value = 0  # H155
$ws.Cells.Item(156, 2).Value = 0  # B156
$ws.Cells.Item(156, 5).Value = 0  # E156
$ws.Cells.Item(157, 2).Value = 0  # B157
$ws.Cells.Item(157, 5).Value = 0  # E157
$ws.Cells.Item(157, 8).Value = 0  # H157
$ws.Cells.Item(158, 2).Value = 0  # B158
$ws.Cells.Item(158, 5).Value = 0.18991837978  # E158
$ws.Cells.Item(159, 5).Value = 0.3636960887  # E159
$ws.Cells.Item(159, 8).Value = 0  # H159
$ws.Cells.Item(160, 5).Value = 0.7873400355  # E160
$ws.Cells.Item(161, 5).Value = 0.8977813776  # E161
